$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2
Set-TextCell $ws "D2" "42.224.26"
Set-TextCell $ws "E2" "  +0.13%  "

# Row 3
Set-TextCell $ws "D3" "2.293.80"

# Row 4
Set-TextCell $ws "E4" "  +0.02%  "

# Row 5
Set-TextCell $ws "D5" "318.11"
Set-TextCell $ws "E5" "  +1.55%  "

# Row 6
Set-TextCell $ws "D6" "102.45"
Set-TextCell $ws "E6" "  -3.87%  "

# Row 8
Set-TextCell $ws "E8" "  -0.09%  "

# Row 9
Set-TextCell $ws "D9" "0.605"
Set-TextCell $ws "E9" "  -0.93%  "

# Row 10
Set-TextCell $ws "E10" "  -1.79%  "

# Row 11
Set-TextCell $ws "D11" "0.0905"
Set-TextCell $ws "E11" "  -0.92%  "

# Row 12
Set-TextCell $ws "D12" "8.41"
Set-TextCell $ws "E12" "  +1.22%  "

# Row 13
Set-TextCell $ws "E13" "  +0.05%  "

# Row 14
Set-TextCell $ws "D14" "0.955"
Set-TextCell $ws "E14" "  -1.88%  "

# Row 15
Set-TextCell $ws "E15" "  -1.98%  "

# Row 16
Set-TextCell $ws "D16" "2.642.14"
Set-TextCell $ws "E16" "  -0.09%  "

# Row 17
Set-TextCell $ws "D17" "2.288.75"
Set-TextCell $ws "E17" "  -1.17%  "

# Row 18
Set-TextCell $ws "D18" "42.343.54"
Set-TextCell $ws "E18" "  +0.86%  "

# Row 19
Set-TextCell $ws "E19" "  -2.11%  "

# Row 20
Set-TextCell $ws "E20" "  +0.69%  "

# Row 21
Set-TextCell $ws "D21" "12.35"
Set-TextCell $ws "E21" "  +31.78%  "

# Row 22
Set-TextCell $ws "D22" "73.38"
Set-TextCell $ws "E22" "  +0.23%  "

# Row 23
Set-TextCell $ws "E23" "  +2.38%  "

# Row 24
Set-TextCell $ws "D24" "275.45"
Set-TextCell $ws "E24" "  +7.06%  "

# Row 25
Set-TextCell $ws "E25" "  -2.95%  "

# Row 26
Set-TextCell $ws "E26" "  -0.28%  "

# Row 27
Set-TextCell $ws "D27" "10.83"
Set-TextCell $ws "E27" "  -1.81%  "

# Row 28
Set-TextCell $ws "D28" "2.36"
Set-TextCell $ws "E28" "  +1.04%  "

# Row 29
Set-TextCell $ws "D29" "22.71"
Set-TextCell $ws "E29" "  -0.61%  "

# Row 30
Set-TextCell $ws "D30" "37.47"
Set-TextCell $ws "E30" "  +4.92%  "

# Row 31
Set-TextCell $ws "D31" "165.75"
Set-TextCell $ws "E31" "  -0.40%  "

# Row 32
Set-TextCell $ws "D32" "6.03"
Set-TextCell $ws "E32" "  +4.01%  "

# Row 33
Set-TextCell $ws "E33" "  -2.19%  "

# Row 34
Set-TextCell $ws "E34" "  +2.87%  "

# Row 35
Set-TextCell $ws "E35" "  -8.56%  "

# Row 36
Set-TextCell $ws "E36" "  -1.45%  "

# Row 37
Set-TextCell $ws "E37" "  -1.02%  "

# Row 38
Set-TextCell $ws "D38" "0.0362"
Set-TextCell $ws "E38" "  +2.17%  "

# Row 39
Set-TextCell $ws "D39" "3.70"
Set-TextCell $ws "E39" "  +1.83%  "

# Row 40
Set-TextCell $ws "E40" "  -5.27%  "

# Row 41
Set-TextCell $ws "E41" "  -0.90%  "

# Row 42
Set-TextCell $ws "B42" "MultiversX"
Set-TextCell $ws "C42" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws "D42" "69.67"
Set-TextCell $ws "E42" "  -3.38%  "

# Row 43
Set-TextCell $ws "B43" "BitcoinSV"
Set-TextCell $ws "C43" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell $ws "D43" "95.69"
Set-TextCell $ws "E43" "  -2.32%  "

# Row 44
Set-TextCell $ws "E44" "  +0.30%  "

# Row 45
Set-TextCell $ws "E45" "  -1.84%  "

# Row 46
Set-TextCell $ws "D46" "11.95"
Set-TextCell $ws "E46" "  -3.39%  "

# Row 47
Set-TextCell $ws "D47" "111.94"
Set-TextCell $ws "E47" "  -1.10%  "

# Row 48
Set-TextCell $ws "D48" "79.14"
Set-TextCell $ws "E48" "  +3.83%  "

# Row 49
Set-TextCell $ws "D49" "8.94"
Set-TextCell $ws "E49" "  -1.99%  "

# Row 50
Set-TextCell $ws "D50" "5.25"
Set-TextCell $ws "E50" "  -1.29%  "

# Row 51
Set-TextCell $ws "D51" "1.599.55"
Set-TextCell $ws "E51" "  +3.59%  "
